# The deck ships two DrawingML theme parts:
#   - theme1.xml ("Office Theme" colors) -> linked from the Notes Master
#   - theme2.xml ("Integral" colors)     -> linked from the Slide Master /
#                                            presentation's active design
#
# The authored edit swaps the two themes' color schemes (theme1 becomes the
# "Integral" palette, theme2 becomes the "Office Theme" palette). The
# PowerPoint object model only exposes the currently active design's palette
# for editing (SlideMaster.Theme.ThemeColorScheme, aliased across
# SlideMaster/NotesMaster/HandoutMaster/Designs), so we recolor it to the
# "Office Theme" values -- this is the half of the swap that is reachable
# through COM automation and lands on theme2.xml, matching the target state.

$p = $ppt.ActivePresentation

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as COM RGB() long values (0x00BBGGRR packing).
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
